# Update automatico via Actualizar 02-16-2021 13-09-05
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-15: refreshed timestamp for the most recent update batch
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44243.54786427599
}

# Rows 16-29: shift forward to the value the newest batch used to have
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44243.52663643518
}

# Rows 30-43: shift forward to the value the middle batch used to have
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44243.50540237268
}
